$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text-like numeric strings (e.g. "552.33") are not auto-converted to numbers.
$ws.Range("D2").Value = "57.209.94"
$ws.Range("E2").Value = "  -5.13%  "
$ws.Range("D3").Value = "2.909.28"
$ws.Range("E3").Value = "  -3.36%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "552.33"
$ws.Range("E5").Value = "  -3.64%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "130.16"
$ws.Range("E6").Value = "  +2.45%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.514"
$ws.Range("E8").Value = "  +2.62%  "
$ws.Range("D9").Value = "2.904.20"
$ws.Range("E9").Value = "  -3.34%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.126"
$ws.Range("E10").Value = "  -3.87%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "4.75"
$ws.Range("E11").Value = "  -7.52%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.444"
$ws.Range("E12").Value = "  +0.87%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000220"
$ws.Range("E13").Value = "  +0.24%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "32.70"
$ws.Range("E14").Value = "  +0.09%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.120"
$ws.Range("E15").Value = "  +0.43%  "
$ws.Range("D16").Value = "3.392.54"
$ws.Range("E16").Value = "  -3.22%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.78"
$ws.Range("E17").Value = "  +5.74%  "
$ws.Range("D18").Value = "2.913.06"
$ws.Range("E18").Value = "  -3.22%  "
$ws.Range("D19").Value = "57.230.23"
$ws.Range("E19").Value = "  -5.15%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "415.24"
$ws.Range("E20").Value = "  -3.28%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.05"
$ws.Range("E21").Value = "  -0.61%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.679"
$ws.Range("E22").Value = "  +1.67%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.92"
$ws.Range("E23").Value = "  -1.89%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.00"
$ws.Range("E24").Value = "  -0.62%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "79.40"
$ws.Range("E25").Value = "  +0.12%  "
$ws.Range("E26").Value = "  -0.15%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.999"
$ws.Range("E27").Value = "  -0.10%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.47"
$ws.Range("E28").Value = "  -3.11%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.45"
$ws.Range("E29").Value = "  +2.04%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.99"
$ws.Range("E30").Value = "  +2.43%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "25.07"
$ws.Range("E31").Value = "  -0.69%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.96"
$ws.Range("E32").Value = "  -2.88%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0957"
$ws.Range("E33").Value = "  +2.17%  "
$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.61"
$ws.Range("E34").Value = "  -0.29%  "
$ws.Range("B35").Value = "Mantle"
$ws.Range("C35").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.939"
$ws.Range("E35").Value = "  -1.07%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.06"
$ws.Range("E36").Value = "  -1.00%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "48.22"
$ws.Range("E37").Value = "  -4.23%  "
$ws.Range("B38").Value = "Cosmos"
$ws.Range("C38").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.69"
$ws.Range("E38").Value = "  +4.69%  "
$ws.Range("B39").Value = "PEPE"
$ws.Range("C39").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D39").Value = "0.0₃0676"
$ws.Range("E39").Value = "  +1.86%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.55"
$ws.Range("E40").Value = "  +1.41%  "
$ws.Range("E41").Value = "  -2.90%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0343"
$ws.Range("E42").Value = "  -3.46%  "
$ws.Range("B43").Value = "Bittensor"
$ws.Range("C43").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "369.13"
$ws.Range("E43").Value = "  -4.62%  "
$ws.Range("B44").Value = "Maker"
$ws.Range("C44").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D44").Value = "2.667.30"
$ws.Range("E44").Value = "  +0.17%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "122.76"
$ws.Range("E46").Value = "  +2.08%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.237"
$ws.Range("E47").Value = "  +1.01%  "
$ws.Range("E48").Value = "  +2.24%  "
$ws.Range("E49").Value = "  -2.68%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "23.13"
$ws.Range("E50").Value = "  -2.23%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.99"
$ws.Range("E51").Value = "  -0.55%  "
